# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" sheet (between "2021-Q2" and "总计") with the
#   quarter's fund-holding detail rows.
# - Refresh the "总计" (totals) summary sheet with a new leading row for
#   2022-Q1 (count=7, value=0.25 亿元), pushing the older rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value as TEXT (keeps values like "10.39" / "000968"
# from being silently reinterpreted as numbers and losing their exact
# string form / leading zeros).
# ---------------------------------------------------------------------
function Set-TextValue($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Helper: apply the "header / index column" look used throughout this
# workbook - bold font, thin box border, centered + top aligned. Applied
# one cell at a time (rather than to a multi-cell range in one shot) so
# every cell gets the identical, fully-boxed style instead of a mix of
# "shared interior edge" border variants.
function Set-HeaderStyle($rng) {
    foreach ($cell in $rng.Cells) {
        $cell.Font.Bold = $true
        $cell.HorizontalAlignment = -4108
        $cell.VerticalAlignment = -4160
        $cell.Borders.Item(7).LineStyle = 1
        $cell.Borders.Item(8).LineStyle = 1
        $cell.Borders.Item(9).LineStyle = 1
        $cell.Borders.Item(10).LineStyle = 1
    }
}

# ---------------------------------------------------------------------
# 1. Drop the existing "总计" sheet - it gets rebuilt (with the new
#    2022-Q1 row prepended) after the 2022-Q1 sheet so the sheet order /
#    sheetId allocation matches.
# ---------------------------------------------------------------------
$totalOld = $wb.Worksheets.Item("总计")
$totalOld.Delete()

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q1" sheet right after "2021-Q2".
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1 = $wb.Worksheets.Add($null, $afterSheet)
$q1.Name = "2022-Q1"

$q1Header = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $q1Header) {
    $q1.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}
Set-HeaderStyle $q1.Range("B1:H1")

$q1Rows = @(
    @("000968", "广发中证养老产业指数A", "10.39", "94.08", "1.48", "0.1538", 7),
    @("014053", "太平睿庆混合A", "2.46", "31.00", "1.30", "0.0320", 5),
    @("010690", "万家互联互通核心资产量化策略混合A", "0.85", "94.05", "3.02", "0.0257", 8),
    @("002982", "广发中证养老产业指数C", "0.88", "94.08", "1.48", "0.0130", 7),
    @("516560", "华宝养老ETF", "0.75", "97.92", "1.54", "0.0116", 7),
    @("014054", "太平睿庆混合C", "0.67", "31.00", "1.30", "0.0087", 5),
    @("010691", "万家互联互通核心资产量化策略混合C", "0.20", "94.05", "3.02", "0.0060", 8)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Cells.Item($r, 1).Value = $r - 2
    Set-TextValue $q1 $r 2 $row[0]
    Set-TextValue $q1 $r 3 $row[1]
    Set-TextValue $q1 $r 4 $row[2]
    Set-TextValue $q1 $r 5 $row[3]
    Set-TextValue $q1 $r 6 $row[4]
    Set-TextValue $q1 $r 7 $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}
Set-HeaderStyle $q1.Range("A2:A8")

# ---------------------------------------------------------------------
# 3. Rebuild the "总计" sheet after "2022-Q1", with the new 2022-Q1 row
#    on top followed by the previously existing rows.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$totalHeader = @("日期", "持有数量(只)", "持有市值(亿元)")
$col = 2
foreach ($h in $totalHeader) {
    $total.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}
Set-HeaderStyle $total.Range("B1:D1")

$totalRows = @(
    @("2022-Q1", 7, 0.25),
    @("2021-Q2", 1, 0),
    @("2021-Q1", 6, 0.46),
    @("2020-Q4", 8, 0.59)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $r - 2
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}
Set-HeaderStyle $total.Range("A2:A5")

Write-Output "2022-Q1 sheet added; 总计 refreshed"
